# "multiplicadors a medi modificats" - update the "terciari" (column G)
# treatment-multiplier codes for a handful of EDAR rows, and correct the
# "poblacio_sanejada" figure for Sant Quirze de Besora (row 61, D61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose tertiary-treatment combo (column G) is removed entirely.
$ws.Range("G13").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("G25").ClearContents()
$ws.Range("G34").ClearContents()
$ws.Range("G38").ClearContents()

# Rows whose tertiary-treatment combo (column G) changes value.
$ws.Range("G33").Value = "UF,UV"
$ws.Range("G49").Value = "UF,UV"
$ws.Range("G57").Value = "O3,SF"
$ws.Range("G70").Value = "SF,UV"
$ws.Range("G78").Value = "UF,UV"

# Population served ("poblacio_sanejada") correction for Sant Quirze de Besora.
$ws.Range("D61").Value = 3815
